$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sample order data (rows 2-16) with generic placeholder
# field-name values, matching the new "database access" column mapping.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = "sku"
    $ws.Cells.Item($r, 2).Value = "name"
    $ws.Cells.Item($r, 3).Value = "quantity"
    $ws.Cells.Item($r, 4).Value = "cost_per"
    $ws.Cells.Item($r, 5).Value = "total_cost"
}

# Remove the now-superfluous 17th row entirely (shrinks used range to E16).
$ws.Rows.Item(17).Delete()
